# Apply scheduled-runner price/profit updates to Behemoth_Profits sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Range("H6").Value = 861.3
$ws.Range("I6").Value = 870.6667
$ws.Range("J6").Value = 777
$ws.Range("K6").Value = 2612.0001
$ws.Range("L6").Value = 2331
$ws.Range("M6").Value = -2500.0001
$ws.Range("N6").Value = -2555

# row 8
$ws.Range("H8").Value = 967.4783
$ws.Range("I8").Value = 967.4783
$ws.Range("K8").Value = 2902.4349
$ws.Range("M8").Value = -2763.4349

# row 45
$ws.Range("H45").Value = 2000
$ws.Range("J45").Value = 2000
$ws.Range("L45").Value = 6000
$ws.Range("N45").Value = -6384

# row 69
$ws.Range("H69").Value = 14594.857
$ws.Range("I69").Value = 3600
$ws.Range("J69").Value = 29254.666
$ws.Range("K69").Value = 10800
$ws.Range("L69").Value = 87763.99800000001
$ws.Range("M69").Value = -9926
$ws.Range("N69").Value = -89511.99800000001

# row 72
$ws.Range("H72").Value = 14594.857
$ws.Range("I72").Value = 3600
$ws.Range("J72").Value = 29254.666
$ws.Range("K72").Value = 32400
$ws.Range("L72").Value = 263291.994
$ws.Range("M72").Value = -28032
$ws.Range("N72").Value = -272027.994

# row 75
$ws.Range("H75").Value = 48578.5
$ws.Range("J75").Value = 48578.5
$ws.Range("L75").Value = 48578.5
$ws.Range("N75").Value = -50450.5

# row 78
$ws.Range("H78").Value = 48578.5
$ws.Range("J78").Value = 48578.5
$ws.Range("L78").Value = 145735.5
$ws.Range("N78").Value = -155095.5

# row 100
$ws.Range("H100").Value = 3417
$ws.Range("I100").Value = 1002.5
$ws.Range("J100").Value = 3899.9
$ws.Range("K100").Value = 1002.5
$ws.Range("L100").Value = 3899.9
$ws.Range("M100").Value = -461.5
$ws.Range("N100").Value = -4981.9

# row 112
$ws.Range("H112").Value = 1701.4857
$ws.Range("J112").Value = 1824.4333
$ws.Range("L112").Value = 5473.2999
$ws.Range("N112").Value = -7689.2999

# row 132
$ws.Range("H132").Value = 2258.5
$ws.Range("I132").Value = 2260.6592
$ws.Range("K132").Value = 6781.9776
$ws.Range("M132").Value = -4251.9776

# row 137
$ws.Range("H137").Value = 5453.48
$ws.Range("I137").Value = 2792.9443
$ws.Range("K137").Value = 8378.832900000001
$ws.Range("M137").Value = -5828.832900000001

# row 138
$ws.Range("H138").Value = 1997.72
$ws.Range("J138").Value = 2658.8035
$ws.Range("L138").Value = 7976.4105
$ws.Range("N138").Value = -18256.4105


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 8077945
$ws.Range("I32").Value = 8941682
$ws.Range("J32").Value = 16397.5
$ws.Range("K32").Value = 8941682
$ws.Range("L32").Value = 16397.5
$ws.Range("M32").Value = -8941395
$ws.Range("N32").Value = -16971.5

# row 61
$ws.Range("H61").Value = 9457145
$ws.Range("I61").Value = 11908379
$ws.Range("K61").Value = 11908379
$ws.Range("M61").Value = -11908167

# row 74
$ws.Range("H74").Value = 5421473
$ws.Range("I74").Value = 6759116.5
$ws.Range("J74").Value = 922127.2
$ws.Range("K74").Value = 6759116.5
$ws.Range("L74").Value = 922127.2
$ws.Range("M74").Value = -6758242.5
$ws.Range("N74").Value = -923875.2

# row 77
$ws.Range("H77").Value = 5421473
$ws.Range("I77").Value = 6759116.5
$ws.Range("J77").Value = 922127.2
$ws.Range("K77").Value = 33795582.5
$ws.Range("L77").Value = 4610636
$ws.Range("M77").Value = -33791214.5
$ws.Range("N77").Value = -4619372

# row 102
$ws.Range("H102").Value = 7435.231
$ws.Range("I102").Value = 9027.666999999999
$ws.Range("J102").Value = 747
$ws.Range("K102").Value = 9027.666999999999
$ws.Range("L102").Value = 747
$ws.Range("M102").Value = -7405.666999999999
$ws.Range("N102").Value = -3991

# row 122
$ws.Range("H122").Value = 2295.7334
$ws.Range("I122").Value = 2295.7334
$ws.Range("K122").Value = 6887.2002
$ws.Range("M122").Value = -4437.2002

# row 136
$ws.Range("H136").Value = 9457145
$ws.Range("I136").Value = 11908379
$ws.Range("K136").Value = 35725137
$ws.Range("M136").Value = -35722587


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 27
$ws.Range("H27").Value = 72494.5
$ws.Range("J27").Value = 72494.5
$ws.Range("L27").Value = 72494.5
$ws.Range("N27").Value = -72878.5

# row 64
$ws.Range("H64").Value = 2288.353
$ws.Range("J64").Value = 1614.125
$ws.Range("L64").Value = 1614.125
$ws.Range("N64").Value = -2064.125

# row 67
$ws.Range("H67").Value = 2288.353
$ws.Range("J67").Value = 1614.125
$ws.Range("L67").Value = 1614.125
$ws.Range("N67").Value = -3174.125

# row 94
$ws.Range("H94").Value = 1833.9231
$ws.Range("I94").Value = 2140.6667
$ws.Range("K94").Value = 2140.6667
$ws.Range("M94").Value = -1689.6667


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 517206.03
$ws.Range("I31").Value = 12232.417
$ws.Range("K31").Value = 12232.417
$ws.Range("M31").Value = -11937.417

# row 34
$ws.Range("H34").Value = 517206.03
$ws.Range("I34").Value = 12232.417
$ws.Range("K34").Value = 12232.417
$ws.Range("M34").Value = -12030.417

# row 132
$ws.Range("H132").Value = 2050.919
$ws.Range("I132").Value = 2021.742
$ws.Range("K132").Value = 6065.226
$ws.Range("M132").Value = -3535.226

# row 134
$ws.Range("H134").Value = 2664.432
$ws.Range("I134").Value = 1855.4517
$ws.Range("K134").Value = 5566.355100000001
$ws.Range("M134").Value = -3031.355100000001


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 4539.222
$ws.Range("J131").Value = 3866.309
$ws.Range("L131").Value = 11598.927
$ws.Range("N131").Value = -21678.927

# row 137
$ws.Range("H137").Value = 4783.1113
$ws.Range("I137").Value = 4006.25
$ws.Range("K137").Value = 12018.75
$ws.Range("M137").Value = -6918.75

# row 139
$ws.Range("H139").Value = 2161.484
$ws.Range("I139").Value = 1533.619
$ws.Range("J139").Value = 3480
$ws.Range("K139").Value = 4600.857
$ws.Range("L139").Value = 10440
$ws.Range("M139").Value = 539.143
$ws.Range("N139").Value = -20720

# row 140
$ws.Range("H140").Value = 87548.77
$ws.Range("I140").Value = 87548.77
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 262646.31
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -257466.31
$ws.Range("N140").ClearContents()


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5254.375
$ws.Range("I70").Value = 5005.8335
$ws.Range("K70").Value = 5005.8335
$ws.Range("M70").Value = -4735.8335

# row 73
$ws.Range("H73").Value = 5254.375
$ws.Range("I73").Value = 5005.8335
$ws.Range("K73").Value = 5005.8335
$ws.Range("M73").Value = -4069.8335


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 48170.17
$ws.Range("I136").Value = 5457.2
$ws.Range("J136").Value = 315126.25
$ws.Range("K136").Value = 16371.6
$ws.Range("L136").Value = 945378.75
$ws.Range("M136").Value = -13821.6
$ws.Range("N136").Value = -950478.75


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 15
$ws.Range("H15").Value = 55003.5
$ws.Range("J15").Value = 55003.5
$ws.Range("L15").Value = 55003.5
$ws.Range("N15").Value = -55579.5

# row 100
$ws.Range("H100").Value = 1479.1765
$ws.Range("I100").Value = 1776.909
$ws.Range("K100").Value = 3553.818
$ws.Range("M100").Value = -3012.818

# row 132
$ws.Range("H132").Value = 3848.2273
$ws.Range("I132").Value = 2447.3572
$ws.Range("K132").Value = 7342.071599999999
$ws.Range("M132").Value = -4812.071599999999

# row 136
$ws.Range("H136").Value = 2808.3684
$ws.Range("I136").Value = 2848.5862
$ws.Range("K136").Value = 8545.758600000001
$ws.Range("M136").Value = -5995.758600000001
